# Apply updated survey figures to rows 9, 10, 11, 16 and relabel row 18.
# All values in this sheet are stored as text (inline strings), so we
# prefix numeric-looking values with a leading apostrophe to force Excel
# to keep them as text rather than silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - Primary
$ws.Range("B9").Value  = "'0.3"
$ws.Range("C9").Value  = "'0.22"
$ws.Range("D9").Value  = "'0.28"
$ws.Range("E9").Value  = "'0.23"
$ws.Range("F9").Value  = "'0.24"
$ws.Range("G9").Value  = "'0.12"
$ws.Range("H9").Value  = "'0.1"
$ws.Range("I9").Value  = "'0.15"
$ws.Range("J9").Value  = "'0.14"
$ws.Range("K9").Value  = "'0.06"
$ws.Range("L9").Value  = "'0.08"

# Row 10 - Secondary
$ws.Range("B10").Value = "'0.22"
$ws.Range("C10").Value = "'0.43"
$ws.Range("D10").Value = "'0.33"
$ws.Range("E10").Value = "'0.44"
$ws.Range("F10").Value = "'0.35"
$ws.Range("G10").Value = "'0.51"
$ws.Range("I10").Value = "'0.57"
$ws.Range("J10").Value = "'0.61"
$ws.Range("L10").Value = "'0.56"

# Row 11 - University
$ws.Range("B11").Value = "'0.05"
$ws.Range("C11").Value = "'0.15"
$ws.Range("D11").Value = "'0.09"
$ws.Range("E11").Value = "'0.18"
$ws.Range("F11").Value = "'0.1"
$ws.Range("G11").Value = "'0.25"
$ws.Range("H11").Value = "'0.35"
$ws.Range("I11").Value = "'0.16"
$ws.Range("J11").Value = "'0.18"
$ws.Range("K11").Value = "'0.31"
$ws.Range("L11").Value = "'0.29"

# Row 16 - Married/Cohabiting
$ws.Range("B16").Value = "'0.59"
$ws.Range("C16").Value = "'0.44"
$ws.Range("D16").Value = "'0.47"
$ws.Range("E16").Value = "'0.45"
$ws.Range("F16").Value = "'0.5"
$ws.Range("G16").Value = "'0.55"
$ws.Range("H16").Value = "'0.63"
$ws.Range("I16").Value = "'0.52"
$ws.Range("J16").Value = "'0.37"
$ws.Range("K16").Value = "'0.61"
$ws.Range("L16").Value = "'0.51"

# Row 18 label - "24 - 49" corrected to "25 - 49"
$ws.Range("A18").Value = "25 - 49"
